$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.237.80"
$ws.Range("E2").Value = "  +1.87%  "
$ws.Range("D3").Value = "1.895.01"
$ws.Range("E3").Value = "  -0.61%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.15%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5177"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4005"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08394"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.114"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +11.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.443"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.11%  "
$ws.Range("D14").Value = "1.893.91"
$ws.Range("E14").Value = "  -0.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.324"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06639"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.48%  "
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.952"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.27%  "
$ws.Range("D23").Value = "30.230.97"
$ws.Range("E23").Value = "  +1.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.227"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("D26").Value = "2.112.19"
$ws.Range("E26").Value = "  -0.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.352"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.091"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.89%  "
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.084"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.747"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02493"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06540"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.285"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2193"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.31%  "
$ws.Range("E39").Value = "  -0.96%  "
$ws.Range("E40").Value = "  +4.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.732"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6499"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.230"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6098"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.700"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.64%  "
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.236"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.163"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.20%  "
